$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "28.499.03"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.875.26"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.52%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.009"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.80%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "315.92"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("E6").Value = "  -0.80%  "

$ws.Range("E7").Value = "  -0.59%  "

$ws.Range("E8").Value = "  -0.41%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.08364"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("E10").Value = "  -1.07%  "

$ws.Range("E11").Value = "  -0.77%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "6.225"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.19%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.874.93"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.08%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "20.41"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

$ws.Range("E15").Value = "  +0.79%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.009"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("E18").Value = "  +0.14%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06731"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.17%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.71"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.41%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.79%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.925"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "28.512.67"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.10"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("E25").Value = "  -1.19%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.088.87"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.20%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "162.03"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("E28").Value = "  +0.33%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.387"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -0.76%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "125.77"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.78%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.1044"
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.040"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.768"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.54%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.614"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("E35").Value = "  +0.74%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.06544"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.2160"
$cell.Style = "Normal"

$ws.Range("E38").Value = "  -4.21%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "5.052"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +1.83%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.191"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.245"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6394"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  -0.22%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.6007"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "13.05"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.690"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("E48").Value = "  +1.05%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.214"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "121.82"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "

$ws.Range("E51").Value = "  -11.00%  "
